$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at sheet row 148, shifting rows 148:214 down to 149:215.
$ws.Range("A148").EntireRow.Insert()

# Populate the newly inserted row 148 with its data (same template as the
# surrounding "Ciboulette" rows, but with its own Fecha/Volumen values).
$ws.Range("A148").Value = 3
$ws.Range("B148").Value = "Femacal de La Calera"
$ws.Range("C148").Value = "Coquimbo"
$ws.Range("D148").Value = 44523
$ws.Range("D148").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E148").Value = 5
$ws.Range("F148").Value = 100112039
$ws.Range("G148").Value = "Ciboulette"
$ws.Range("H148").Value = "Sin especificar"
$ws.Range("I148").Value = "Primera"
$ws.Range("J148").Value = 160
$ws.Range("K148").Value = 1500
$ws.Range("L148").Value = 1500
$ws.Range("M148").Value = 1500
$ws.Range("N148").Value = "`$/docena de atados"
$ws.Range("O148").Value = "Provincia de Quillota"
$ws.Range("P148").Value = 500
$ws.Range("Q148").Value = 3
$ws.Range("R148").Value = "Hortaliza"
